# Updates the "Poker - Year Figures" workbook's Sheet1 data for the 2025
# year-to-date block (rows 201-210) to reflect the refreshed Power Query
# figures pulled for the "12 - Dec" month-end snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 201 - Andy
$ws.Range("D201").Value = 53
$ws.Range("F201").Value = 53
$ws.Range("G201").Value = 170200
$ws.Range("H201").Value = 230

# Row 202 - Prashant
$ws.Range("D202").Value = 52
$ws.Range("F202").Value = 52
$ws.Range("G202").Value = 153400
$ws.Range("H202").Value = 160

# Row 203 - Richard
$ws.Range("D203").Value = 48
$ws.Range("F203").Value = 48
$ws.Range("G203").Value = 147950
$ws.Range("H203").Value = 195
$ws.Range("I203").Value = 95

# Row 204 - Pepe
$ws.Range("D204").Value = 37
$ws.Range("F204").Value = 37
$ws.Range("G204").Value = 115400
$ws.Range("I204").Value = -40

# Row 205 - Maisy
$ws.Range("D205").Value = 35
$ws.Range("F205").Value = 35
$ws.Range("G205").Value = 101550
$ws.Range("I205").Value = -30

# Row 206 - Mark (SRank moves from 6 to 5, tied with Maisy)
$ws.Range("C206").Value = 5
$ws.Range("D206").Value = 35
$ws.Range("F206").Value = 35
$ws.Range("G206").Value = 101550
$ws.Range("H206").Value = 85
$ws.Range("I206").Value = -15

# Row 207 - was Matt, now Anthony (swapped with row 208)
$ws.Range("B207").Value = "Anthony"
$ws.Range("D207").Value = 32
$ws.Range("F207").Value = 32
$ws.Range("G207").Value = 81450
$ws.Range("H207").Value = 30
$ws.Range("I207").Value = -50
$ws.Range("K207").Value = 350

# Row 208 - was Anthony, now Matt (swapped with row 207)
$ws.Range("B208").Value = "Matt"
$ws.Range("D208").Value = 30
$ws.Range("F208").Value = 30
$ws.Range("G208").Value = 89750
$ws.Range("H208").Value = 60
$ws.Range("I208").Value = -30
$ws.Range("K208").Value = 362

# Row 209 - Jon
$ws.Range("D209").Value = 19
$ws.Range("F209").Value = 19
$ws.Range("G209").Value = 60550
$ws.Range("I209").Value = -80

# Row 210 - Alex
$ws.Range("D210").Value = 18
$ws.Range("F210").Value = 18
$ws.Range("G210").Value = 64650
$ws.Range("I210").Value = -30
